$d = $word.ActiveDocument
$d.Content.Find.Execute("UI enemy states and detection bar:", $true, $false, $false, $false, $false,
                         $true, 1, $false, "UI enemy states and detection:", 2)
